# Clarify cascading select example
# Rewrites the survey/choices/settings sheets of the cascading-select demo
# workbook so that:
#  - the survey sheet uses clearer field names/labels and filters based on
#    the newly-named fields
#  - the choices sheet splits the combined "choice_filter" column into
#    separate "state" and "county" columns
#  - the settings sheet gets a friendlier title and an explicit form_id

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# survey sheet
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# row 2: select_one states
$survey.Range("B2").Value = "selected_state"
$survey.Range("C2").Value = "Select a state"

# row 3: select_one counties
$survey.Range("B3").Value = "selected_county"
$survey.Range("C3").Value = "Select a county"
$survey.Range("D3").Value = "state=`${selected_state}"

# row 4: select_one cities
$survey.Range("B4").Value = "selected_city"
$survey.Range("C4").Value = "Select a city"
$survey.Range("D4").Value = "county=`${selected_county}"

# ---------------------------------------------------------------------
# choices sheet
# ---------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

# header row: choice_filter -> state, and add a new county column
$choices.Range("D1").Value = "state"

$choices.Range("D1").Copy()
$choices.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$choices.Range("E1").Value = "county"

# states list (rows 2-3) is unchanged

# counties list: drop the old row 4 (king/washington), shift rows up,
# and add brewster/cameron at the bottom
$choices.Range("A4:D4").ClearContents()

$choices.Range("B5").Value = "king"
$choices.Range("C5").Value = "King"
$choices.Range("D5").Value = "washington"

$choices.Range("B6").Value = "pierce"
$choices.Range("C6").Value = "Pierce"
$choices.Range("D6").Value = "washington"

$choices.Range("B7").Value = "brewster"
$choices.Range("C7").Value = "Brewster"
$choices.Range("D7").Value = "texas"

$choices.Range("A8").Value = "counties"
$choices.Range("B8").Value = "cameron"
$choices.Range("C8").Value = "Cameron"
$choices.Range("D8").Value = "texas"

# cities list: drop dumont/finney (old rows 8-9), shift up, and add the
# new "county" column value alongside the existing "state" value
$choices.Range("A9:D9").ClearContents()

$choices.Range("A10").Value = "cities"
$choices.Range("D10").Value = "texas"
$choices.Range("D10").Copy()
$choices.Range("E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$choices.Range("E10").Value = "cameron"

$choices.Range("A11").Value = "cities"
$choices.Range("D11").Value = "texas"
$choices.Range("D11").Copy()
$choices.Range("E11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$choices.Range("E11").Value = "cameron"

$choices.Range("A12").Value = "cities"
$choices.Range("D12").Value = "washington"
$choices.Range("D12").Copy()
$choices.Range("E12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$choices.Range("E12").Value = "king"

$choices.Range("A13").Value = "cities"
$choices.Range("D13").Value = "washington"
$choices.Range("D13").Copy()
$choices.Range("E13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$choices.Range("E13").Value = "king"

$choices.Range("A14").Value = "cities"
$choices.Range("D14").Value = "washington"
$choices.Range("D14").Copy()
$choices.Range("E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$choices.Range("E14").Value = "pierce"

$choices.Range("A15").Value = "cities"
$choices.Range("D15").Value = "washington"
$choices.Range("D15").Copy()
$choices.Range("E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$choices.Range("E15").Value = "pierce"

# ---------------------------------------------------------------------
# settings sheet
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")

$settings.Range("A2").Value = "Cascading select example"

$settings.Range("A2").Copy()
$settings.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$settings.Range("B2").Value = "cascading_select"
